# Junction_Flooding_409 edit: refresh simulation dataset (custom accuracy + new 1000-point run)
# - replace rows 2-5 with newly-simulated values (new start date, finer precision)
# - drop row 6 (dataset now spans 4 timesteps instead of 5)
# - widen most data columns from 7 to 8 chars (col T: 8 -> 9) to fit new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 6); Excel will shift the dimension/used range down automatically
$ws.Rows("6:6").Delete()

# New values for rows 2-5, columns A (Time) through AH (J33)
$data = New-Object "object[,]" 4,34
# row 2
$data[0,0] = 45090.50694444445
$data[0,1] = 6.427
$data[0,2] = 5.211
$data[0,3] = 0
$data[0,4] = 14.561
$data[0,5] = 11.533
$data[0,6] = 4.486
$data[0,7] = 15.002
$data[0,8] = 8.211
$data[0,9] = 4.17
$data[0,10] = 5.579
$data[0,11] = 6.239
$data[0,12] = 6.747
$data[0,13] = 1.607
$data[0,14] = 5.615
$data[0,15] = 7.139
$data[0,16] = 4.876
$data[0,17] = 0.271
$data[0,18] = 0.709
$data[0,19] = 78.413
$data[0,20] = 15.366
$data[0,21] = 5.184
$data[0,22] = 9.362
$data[0,23] = 6.286
$data[0,24] = 0.958
$data[0,25] = 8.106999999999999
$data[0,26] = 4.01
$data[0,27] = 5.031
$data[0,28] = 6.902
$data[0,29] = 7.612
$data[0,30] = 1.257
$data[0,31] = 12.304
$data[0,32] = 3.802
$data[0,33] = 5.737
# row 3
$data[1,0] = 45090.51388888889
$data[1,1] = 22.129
$data[1,2] = 16.792
$data[1,3] = 0.552
$data[1,4] = 48.665
$data[1,5] = 39.798
$data[1,6] = 17.038
$data[1,7] = 63.85
$data[1,8] = 27.051
$data[1,9] = 12.642
$data[1,10] = 18.152
$data[1,11] = 19.686
$data[1,12] = 20.884
$data[1,13] = 5.598
$data[1,14] = 17.647
$data[1,15] = 24.792
$data[1,16] = 14.833
$data[1,17] = 0.245
$data[1,18] = 0.954
$data[1,19] = 261.636
$data[1,20] = 49.332
$data[1,21] = 16.29
$data[1,22] = 32.87
$data[1,23] = 17.946
$data[1,24] = 2.504
$data[1,25] = 32.445
$data[1,26] = 14.185
$data[1,27] = 13.298
$data[1,28] = 15.903
$data[1,29] = 21.158
$data[1,30] = 0.647
$data[1,31] = 57.861
$data[1,32] = 9.675000000000001
$data[1,33] = 20.034
# row 4
$data[2,0] = 45090.52083333334
$data[2,1] = 13.126
$data[2,2] = 9.967000000000001
$data[2,3] = 0.291
$data[2,4] = 28.985
$data[2,5] = 23.621
$data[2,6] = 10.054
$data[2,7] = 43.646
$data[2,8] = 16.086
$data[2,9] = 7.628
$data[2,10] = 10.779
$data[2,11] = 11.728
$data[2,12] = 12.439
$data[2,13] = 3.342
$data[2,14] = 10.504
$data[2,15] = 14.751
$data[2,16] = 8.877000000000001
$data[2,17] = 0.137
$data[2,18] = 0.598
$data[2,19] = 152.92
$data[2,20] = 29.525
$data[2,21] = 9.696999999999999
$data[2,22] = 19.612
$data[2,23] = 10.728
$data[2,24] = 1.544
$data[2,25] = 21.178
$data[2,26] = 8.462
$data[2,27] = 8.023
$data[2,28] = 9.538
$data[2,29] = 12.617
$data[2,30] = 0.435
$data[2,31] = 39.701
$data[2,32] = 5.781
$data[2,33] = 11.916
# row 5
$data[3,0] = 45090.52777777778
$data[3,1] = 20.4
$data[3,2] = 15.4
$data[3,3] = 0.5600000000000001
$data[3,4] = 44.74
$data[3,5] = 36.71
$data[3,6] = 15.84
$data[3,7] = 62.28
$data[3,8] = 24.86
$data[3,9] = 11.49
$data[3,10] = 16.66
$data[3,11] = 18.02
$data[3,12] = 19.09
$data[3,13] = 5.17
$data[3,14] = 16.14
$data[3,15] = 22.88
$data[3,16] = 13.52
$data[3,17] = 0.19
$data[3,18] = 0.76
$data[3,19] = 238.75
$data[3,20] = 45.19
$data[3,21] = 14.9
$data[3,22] = 30.32
$data[3,23] = 16.26
$data[3,24] = 2.27
$data[3,25] = 30.61
$data[3,26] = 13.1
$data[3,27] = 11.96
$data[3,28] = 14.11
$data[3,29] = 19.16
$data[3,30] = 0.33
$data[3,31] = 56.39
$data[3,32] = 8.69
$data[3,33] = 18.48

$ws.Range("A2:AH5").Value2 = $data

# Column width adjustments (stored XML width = ColumnWidth + 0.83 for this workbook/font)
$cols8 = @("B:C","F:G","I:M","O:Q","U:U","W:X","Z:AD","AF:AF","AH:AH")
foreach ($rng in $cols8) {
    $ws.Range($rng).EntireColumn.ColumnWidth = 8 - 0.83
}
$ws.Range("T:T").EntireColumn.ColumnWidth = 9 - 0.83

